$wb = $excel.ActiveWorkbook

# Rename sheets (by index to be robust against any name collisions)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477806120553"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477807999646"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778080005555"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778080615535"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778081245556"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778060845585.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778061035886.csv"
$ws1.Range("B4").Value = "go_stims-1650477806105554.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778061195526.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_2-16504778064565892.csv"
$ws2.Range("B3").Value = "OB-1650477807284552.csv"
$ws2.Range("B4").Value = "TB-16504778079855871.csv"
$ws2.Range("B5").Value = "TB-16504778073605888.csv"
$ws2.Range("B6").Value = "ZB-match_0-16504778064275532.csv"
$ws2.Range("B7").Value = "OB-16504778070745873.csv"
$ws2.Range("B8").Value = "OB-16504778071745527.csv"
$ws2.Range("B9").Value = "ZB-match_1-16504778063115876.csv"
$ws2.Range("B10").Value = "TB-16504778075065882.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778080145555.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778080025582.csv"
$ws4.Range("B4").Value = "MM_stims-1650477808045555.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778080155568.csv"
$ws4.Range("B6").Value = "MM_stims-16504778080605888.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778080465546.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778080935547.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778080645587.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778081085887.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778080775588.csv"
